# SNAPSHOT_v11.1d_2017-05-10
# Adds a new metrics row (row 8, dated 2017-05-10) to Sheet1, mirroring the
# structure/format of the prior row (row 7), updates a handful of the
# computed/measured values for the new snapshot, widens column A to fit the
# new date values, and leaves the selection where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Duplicate row 7's formatting down into the new row 8 -------------------
# Row 7 is the fully-populated template row; copy formats first so the new
# row inherits the same number formats / wrap text / alignment as its
# neighbours, then fill in this snapshot's values explicitly.
$ws.Range("A7:AL7").Copy() | Out-Null
$ws.Range("A8:AL8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New row 8 values ---------------------------------------------------
$ws.Range("A8").Value = 42865
$ws.Range("B8").Value = 229
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 309
$ws.Range("E8").Value = 67
$ws.Range("F8").Value = 3555
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 202
$ws.Range("I8").Value = 216
$ws.Range("J8").Value = 657
$ws.Range("K8").Value = 130
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 225
$ws.Range("N8").Value = 10
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 10
$ws.Range("Q8").Value = 83
$ws.Range("R8").Value = 15
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 7
$ws.Range("U8").Value = 7
$ws.Range("V8").Value = 83
$ws.Range("W8").Value = 19
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 45
$ws.Range("Z8").Value = 240
$ws.Range("AA8").Value = 130
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 229
$ws.Range("AD8").Value = 5
$ws.Range("AE8").Value = 30
$ws.Range("AF8").Value = 85
$ws.Range("AG8").Value = 503
$ws.Range("AH8").Value = 1
$ws.Range("AI8").Value = 1
$ws.Range("AJ8").Value = 0.44776120000000003
$ws.Range("AK8").Value = 0.27508090000000002
$ws.Range("AL8").Value = 0.1414909

# --- Column A now needs to fit the wider set of dates --------------------
$ws.Columns("A:A").AutoFit() | Out-Null

# --- Window / selection bookkeeping --------------------------------------
$win = $wb.Windows.Item(1)
$win.TabRatio = 230
$win.Height = 8280

$ws.Range("G15").Select() | Out-Null
